$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> column letter -> new text value.
# All values are applied as literal text (matching the workbook's existing
# inline-string storage) so that numeric-looking strings such as "10.28"
# or "0.0982" are not reinterpreted by Excel as numbers.
$updates = @{}

$updates[2] = @{ "D" = "58.694.28"; "E" = "  -1.09%  " }
$updates[3] = @{ "D" = "2.634.82"; "E" = "  -0.10%  " }
$updates[4] = @{ "E" = "  -0.05%  " }
$updates[5] = @{ "D" = "516.62"; "E" = "  +0.23%  " }
$updates[6] = @{ "D" = "145.66"; "E" = "  -2.02%  " }
$updates[7] = @{ "E" = "  +0.25%  " }
$updates[8] = @{ "D" = "0.573"; "E" = "  -0.14%  " }
$updates[9] = @{ "D" = "2.643.21"; "E" = "  -0.88%  " }
$updates[10] = @{ "D" = "6.36"; "E" = "  -2.66%  " }
$updates[11] = @{ "E" = "  -1.49%  " }
$updates[12] = @{ "D" = "0.334"; "E" = "  -2.03%  " }
$updates[13] = @{ "E" = "  -0.17%  " }
$updates[14] = @{ "D" = "3.097.81"; "E" = "  +0.03%  " }
$updates[15] = @{ "D" = "58.709.26"; "E" = "  -0.88%  " }
$updates[16] = @{ "D" = "20.80"; "E" = "  -3.00%  " }
$updates[17] = @{ "E" = "  -2.11%  " }
$updates[18] = @{ "D" = "2.645.10"; "E" = "  -0.94%  " }
$updates[19] = @{ "D" = "347.96"; "E" = "  +0.62%  " }
$updates[20] = @{ "D" = "4.45"; "E" = "  -3.72%  " }
$updates[21] = @{ "D" = "10.22"; "E" = "  -3.06%  " }
$updates[22] = @{ "D" = "6.16"; "E" = "  -0.81%  " }
$updates[23] = @{ "D" = "0.998"; "E" = "  -0.01%  " }
$updates[24] = @{ "D" = "61.71"; "E" = "  +0.88%  " }
$updates[25] = @{ "E" = "  -3.15%  " }
$updates[26] = @{ "E" = "  +1.37%  " }
$updates[27] = @{ "D" = "0.996"; "E" = "  +0.09%  " }
$updates[28] = @{ "E" = "  -4.05%  " }
$updates[29] = @{ "D" = "7.01"; "E" = "  -2.20%  " }
$updates[30] = @{ "E" = "  +0.13%  " }
$updates[31] = @{ "E" = "  -5.58%  " }
$updates[32] = @{ "D" = "18.87"; "E" = "  -0.91%  " }
$updates[33] = @{ "E" = "  -0.21%  " }
$updates[34] = @{ "D" = "148.95"; "E" = "  -0.41%  " }
$updates[35] = @{ "D" = "0.989"; "E" = "  -6.42%  " }
$updates[36] = @{ "E" = "  -2.28%  " }
$updates[37] = @{ "E" = "  -1.90%  " }
$updates[38] = @{ "D" = "36.53"; "E" = "  +0.06%  " }
$updates[39] = @{ "D" = "0.837"; "E" = "  -4.80%  " }
$updates[40] = @{ "D" = "1.42"; "E" = "  -1.84%  " }
$updates[41] = @{ "D" = "3.62"; "E" = "  -2.72%  " }
$updates[42] = @{ "D" = "278.64"; "E" = "  -3.28%  " }
$updates[43] = @{ "D" = "0.998"; "E" = "  +0.53%  " }
$updates[44] = @{ "D" = "0.0982" }
$updates[45] = @{ "E" = "  -4.60%  " }
$updates[46] = @{ "D" = "19.57"; "E" = "  -0.84%  " }
$updates[47] = @{ "D" = "0.0523"; "E" = "  -4.25%  " }
$updates[48] = @{ "D" = "10.28"; "E" = "  +0.27%  " }
$updates[49] = @{ "B" = "Maker"; "C" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; "D" = "1.991.31"; "E" = "  +0.43%  " }
$updates[50] = @{ "B" = "VeChain"; "C" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; "D" = "0.0229"; "E" = "  -2.23%  " }
$updates[51] = @{ "E" = "  -3.37%  " }

$colIndex = @{ "B" = 2; "C" = 3; "D" = 4; "E" = 5 }

foreach ($row in $updates.Keys) {
    $rowVals = $updates[$row]
    foreach ($col in $rowVals.Keys) {
        $cell = $ws.Cells.Item($row, $colIndex[$col])
        $cell.NumberFormat = "@"
        $cell.Value = $rowVals[$col]
    }
}
